$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.539.53'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -3.88%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.537.73'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -3.76%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '507.43'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -3.98%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.77'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -7.31%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.564'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -4.17%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.542.30'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -3.91%  '

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -8.45%  '

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -6.20%  '

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.58%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.981.05'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.76%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '58.515.76'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -3.92%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.70'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -5.46%  '

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -5.79%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.540.18'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -3.86%  '

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -4.78%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '337.08'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -4.61%  '

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -5.11%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.998'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.14%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.96'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -4.29%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.49'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.79%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -4.63%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.06%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -4.68%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.650.47'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -3.68%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0787'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -8.94%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.95'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -5.80%  '

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.02%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '149.71'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.44%  '

$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = 'Aptos'
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.83'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -5.11%  '

$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = 'EthereumClassic'
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '18.53'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -4.83%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.54'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -5.09%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.912'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.04%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -5.94%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -7.39%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.07'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.35%  '

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -11.70%  '

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -6.87%  '

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'Bittensor'
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '283.81'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -7.12%  '

$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.53'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -7.47%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0996'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.39%  '

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.03%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.600'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -6.36%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0531'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -5.50%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '18.69'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -5.48%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '10.30'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.40%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0227'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -4.89%  '

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -8.34%  '
